$d = $word.ActiveDocument

# --- 1. Retitle: "IOT-Server Documentation" -> "TechOnIt (Server)" + " Documentation" ---
# Replace the whole title text first (keeps the single run's rPr: sz=32/szCs=32).
$d.Content.Find.Execute("IOT-Server Documentation", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "TechOnIt (Server) Documentation", 2)

# Now split that single run into two runs at the boundary between "(Server)" and
# " Documentation" by nudging formatting on the trailing sub-range (Bold on/off is a
# no-op on the final look but forces Word to materialize a separate run boundary there).
$titlePara = $d.Paragraphs(1).Range
$splitAt = $titlePara.Start + "TechOnIt (Server)".Length
$tail = $d.Range($splitAt, $titlePara.End - 1)
$tail.Bold = 1
$tail.Bold = 0

# --- 2. Add headers/footers (default, even, first) to the lone section (watermark setup) ---
$sec = $d.Sections(1)

# Touching each header/footer range materializes header1/2/3.xml + footer1/2/3.xml
# (even/default/first), their relationships + headerReference/footerReference entries
# in sectPr, and the Header/Footer (+ linked char) styles -- all left empty.
$sec.Headers(1).Range.Text = ""
$sec.Headers(2).Range.Text = ""
$sec.Headers(3).Range.Text = ""
$sec.Footers(1).Range.Text = ""
$sec.Footers(2).Range.Text = ""
$sec.Footers(3).Range.Text = ""
